$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header cells: "_old" columns -> "_FV2410", "_new" columns -> "_FV2504"
$oldHeaders = @(
  "Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old",
  "Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old"
)
$newHeaders = @(
  "Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new",
  "Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
  $col = $i + 1
  $ws.Cells.Item(1, $col).Value = ($oldHeaders[$i] -replace "_old$", "_FV2410")
}
# column 11 is "diff" - untouched
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
  $col = $i + 12
  $ws.Cells.Item(1, $col).Value = ($newHeaders[$i] -replace "_new$", "_FV2504")
}

# 2) Turn the header+data range into an Excel Table ("Table1")
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U93"), $null, 1)
$tbl.TableStyle = ""

# 3) Freeze the header row (pane split inferred from current selection)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
